$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.978.42'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '1.754.95'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.96'
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3823'
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3431'
$ws.Range("E8").Value = '  +0.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.08'
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.121'
$ws.Range("E10").Value = '  -1.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07232'
$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.54'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.170'
$ws.Range("E14").Value = '  -2.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.140'
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").Value = '1.751.27'
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001060'
$ws.Range("E17").Value = '  -1.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06606'
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '79.35'
$ws.Range("E19").Value = '  -3.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.73'
$ws.Range("E21").Value = '  -3.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.210'
$ws.Range("E22").Value = '  -2.86%  '

$ws.Range("D23").Value = '27.992.91'
$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.68'
$ws.Range("E24").Value = '  -3.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.383'
$ws.Range("E25").Value = '  +0.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.73'
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("E27").Value = '  -3.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.307'
$ws.Range("E28").Value = '  -5.37%  '

$ws.Range("D29").Value = '1.951.90'
$ws.Range("E29").Value = '  -0.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.264'
$ws.Range("E30").Value = '  -10.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.40'
$ws.Range("E31").Value = '  -2.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.026'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.847'
$ws.Range("E33").Value = '  -4.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08821'
$ws.Range("E34").Value = '  +0.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.20'
$ws.Range("E35").Value = '  -4.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6569'
$ws.Range("E36").Value = '  -2.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02290'
$ws.Range("E37").Value = '  -4.89%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06171'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.150'
$ws.Range("E39").Value = '  -3.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.513'
$ws.Range("E40").Value = '  +0.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2105'
$ws.Range("E41").Value = '  -3.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.217'
$ws.Range("E42").Value = '  -2.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.942'
$ws.Range("E43").Value = '  -4.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9986'
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.82'
$ws.Range("E45").Value = '  -2.28%  '

$ws.Range("E46").Value = '  +0.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6066'
$ws.Range("E47").Value = '  -2.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.91'
$ws.Range("E48").Value = '  -3.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.010'
$ws.Range("E49").Value = '  -3.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.169'
$ws.Range("E50").Value = '  +2.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.114'
$ws.Range("E51").Value = '  +5.17%  '
